# #4 Update on task division
# Mark the "X" indicators for the pairs of cells (done/not-done style
# checkboxes) that were filled in on rows 8, 10 and 15 of the task
# division sheet, and leave the selection on the cell the author ended
# up on (F9), matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = "X"
$ws.Range("G8").Value = "X"

$ws.Range("M10").Value = "X"
$ws.Range("N10").Value = "X"

$ws.Range("M15").Value = "X"
$ws.Range("N15").Value = "X"

# Restore the view/selection state recorded for the sheet.
$ws.Range("F9").Select()
